$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 77; this shifts rows 77..109 down to 78..110,
# preserving all their existing values (matches the diff's row-shift
# pattern for rows 78-110).
$ws.Rows(77).Insert()

# Populate the newly inserted row 77 with a fresh data record.
$ws.Cells.Item(77, 1).Value = 4
$ws.Cells.Item(77, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(77, 3).Value = "Los Lagos"
$ws.Cells.Item(77, 4).Value = 45089
$ws.Cells.Item(77, 5).Value = 10
$ws.Cells.Item(77, 6).Value = 100112043
$ws.Cells.Item(77, 7).Value = "Pepino dulce"
$ws.Cells.Item(77, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 30
$ws.Cells.Item(77, 11).Value = 17000
$ws.Cells.Item(77, 12).Value = 17000
$ws.Cells.Item(77, 13).Value = 17000
$ws.Cells.Item(77, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(77, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(77, 16).Value = 944
$ws.Cells.Item(77, 17).Value = 18
$ws.Cells.Item(77, 18).Value = "Hortaliza"
